$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right count and Wrong penalty corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right count, Wrong penalty, and summary text corrected
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "52 / 112"
